# feat: incell list, map, and struct
#
# Adds three new in-cell-list/map/struct columns (WarnInfo, IntParams, Tips)
# to the "Activity" sheet, ahead of the existing ClientParam1..3 columns
# (which shift from U:W to X:Z), and renames the "OutputItem" type alias to
# "ExchangeItem" on the "Reward" sheet.

$wb  = $excel.ActiveWorkbook
$act = $wb.Worksheets.Item("Activity")
$rwd = $wb.Worksheets.Item("Reward")

# --- Activity sheet: insert 3 new columns before the old "U" column -------
# This shifts the existing ClientParam1/2/3 columns (U:W) to X:Z, carrying
# their values/styles/column-widths along automatically, matching how the
# workbook was originally edited in Excel.
$act.Columns("U:W").Insert()

# Row 1 - headers
$act.Range("U1").Value = "WarnInfo"
$act.Range("V1").Value = "IntParams"
$act.Range("W1").Value = "Tips"

# Row 2 - field "type" row
$act.Range("U2").Value = "{int32 Id,string Desc,int32 Value}Info"
$act.Range("V2").Value = "[]int32"
$act.Range("W2").Value = "map<int32,string> "

# Row 3 - Chinese description row (U3 intentionally left blank, same as
# source workbook)
$act.Range("V3").Value = "整型参数列表"
$act.Range("W3").Value = "客户端参数1"

# Row 4 - sample data
$act.Range("U4").Value = "1,desc1,100"
$act.Range("V4").Value = 1

# Row 5 - sample data
$act.Range("U5").Value = "2,desc2,500"
$act.Range("V5").Value = 1
$act.Range("W5").Value = "10:hot"

# Row 6 - sample data
$act.Range("U6").Value = "3,desc3,1000"
$act.Range("V6").Value = "1,2,3"
$act.Range("W6").Value = "1:good,2:bad"

# --- Reward sheet: rename the "[OutputItem]int32" type label -------------
$rwd.Range("H2").Value = "[ExchangeItem]int32"

# --- Restore selections / active sheet as left by the edit ---------------
$rwd.Range("M9").Select()
$act.Activate()
$act.Range("U7").Select()
